$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row for Arthur (account 008336728) - worksheet row 8
$ws.Rows.Item(8).Delete()

# After deleting row 8, the Zenilda row (account 004211911) that was on row 12
# shifts up to row 11
$ws.Rows.Item(11).Delete()
